$wb = $excel.ActiveWorkbook

# --- Standup sheet: insert a new weekly column (for 2017-09-28) and fill in
# the feedback / end-of-day status cells ---
$ws1 = $wb.Worksheets.Item("Standup")

[void]$ws1.Columns("C").Insert()

# New date column (matches the weekly cadence of the existing dates)
$ws1.Range("C5").Value = 43006

# "End of Day" status updates / feedback entries
$ws1.Range("C6").Value = "NA"
$ws1.Range("B7").Value = "NA"
$ws1.Range("C7").Value = "End of Day"
$ws1.Range("B8").Value = "NA"
$ws1.Range("C8").Value = "NA"

# Highlight the "End of Day" feedback cells
$ws1.Range("B6").Interior.Color = 255
$ws1.Range("C7").Interior.Color = 255

# --- Professionalism sheet: highlight the "End of Day" feedback cell ---
$ws2 = $wb.Worksheets.Item("Professionalism")
$ws2.Range("E7").Interior.Color = 255

# --- Hardware Development Process sheet: highlight the "End of Day" feedback cells ---
$ws3 = $wb.Worksheets.Item("Hardware Development Process")
$ws3.Range("E7").Interior.Color = 255
$ws3.Range("E8").Interior.Color = 255

# --- Restore/update each sheet's remembered selection, then leave Standup
# as the active (selected) tab, matching the saved view state ---
[void]$ws2.Range("E14").Select()
[void]$ws3.Range("E19").Select()
[void]$ws1.Range("C9").Select()
